$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: merge the two runs in the "extensometro" paragraph into a
# single run with a trailing space, and drop the "_GoBack" bookmark
# that used to sit between them.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBackStart = $goBack.Start

# Delete across the bookmark position (one char before it / one char
# after it == the last letter of "deformacao" + the lone space run).
# This removes both the bookmark and the stray run in one go.
$aroundRange = $d.Range($goBackStart - 1, $goBackStart + 1)
$removedText = $aroundRange.Text
$aroundRange.Delete()

# Put the removed text ("o ") straight back so the paragraph reads
# "...deformação " again, now living inside the single remaining run.
$restoreRange = $d.Range($goBackStart - 1, $goBackStart - 1)
$restoreRange.InsertAfter($removedText)

# ------------------------------------------------------------------
# Change 2: fix the "Ductibilidade" typo inside the results table
# (match case so the lowercase occurrence elsewhere is untouched).
# ------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Execute("Ductibilidade", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Ductilidade", 2) | Out-Null

# ------------------------------------------------------------------
# Change 3: re-create the "_GoBack" bookmark at the end of the
# "Resultados e discussão" heading (last edit location), right after
# the "_Toc33000103" bookmark's end.
# ------------------------------------------------------------------
$tocBookmark = $d.Bookmarks("_Toc33000103")
$headingEnd = $d.Range($tocBookmark.End, $tocBookmark.End)

# Bookmarks.Add("_GoBack", <zero-width range>) is unreliable, so we
# add it around a temporary placeholder character (non-zero width)
# and then delete that character, which collapses the bookmark down
# to zero width while keeping it anchored in the right place.
$headingEnd.InsertAfter("Z")

$placeholderSearch = $d.Content
$placeholderSearch.Find.ClearFormatting()
$placeholderSearch.Find.Execute("discussãoZ", $true, $false, $false, $false, $false, `
    $true, 1, $false, $null, 0) | Out-Null
$placeholderRange = $d.Range($placeholderSearch.End - 1, $placeholderSearch.End)

$d.Bookmarks.Add("_GoBack", $placeholderRange)
$newGoBack = $d.Bookmarks("_GoBack")
$placeholderDelete = $d.Range($newGoBack.Start, $newGoBack.End)
$placeholderDelete.Delete()
